# Refresh crypto price/volume figures on the "cryptos" worksheet.
# Price cells (column D) are forced to Text format before assignment so that
# values which look numeric (e.g. "1.00", "7.80") are stored verbatim as
# strings, matching the source data feed's formatting instead of being
# normalized into numbers by Excel.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '60.081.09'
$ws.Range("E2").Value = '  -3.28%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.286.29'
$ws.Range("E3").Value = '  -4.00%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.00'
$ws.Range("E4").Value = '  -0.01%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '554.75'
$ws.Range("E5").Value = '  -4.09%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '140.47'
$ws.Range("E6").Value = '  -7.85%  '
$ws.Range("E7").Value = '  -0.02%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '3.288.94'
$ws.Range("E8").Value = '  -3.94%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.465'
$ws.Range("E9").Value = '  -3.74%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '7.80'
$ws.Range("E10").Value = '  -3.29%  '
$ws.Range("E11").Value = '  -5.18%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.405'
$ws.Range("E12").Value = '  -2.90%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '3.860.27'
$ws.Range("E13").Value = '  -3.75%  '
$ws.Range("E14").Value = '  -0.13%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '26.79'
$ws.Range("E15").Value = '  -6.71%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '3.297.39'
$ws.Range("E16").Value = '  -3.74%  '
$ws.Range("E17").Value = '  -4.56%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '60.188.25'
$ws.Range("E18").Value = '  -3.22%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '6.07'
$ws.Range("E19").Value = '  -6.69%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '13.84'
$ws.Range("E20").Value = '  -5.07%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '8.53'
$ws.Range("E21").Value = '  -4.58%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '371.85'
$ws.Range("E22").Value = '  -2.96%  '
$ws.Range("E23").Value = '  -1.60%  '
$ws.Range("E24").Value = '  -0.03%  '
$ws.Range("E25").Value = '  -7.17%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '3.435.03'
$ws.Range("E26").Value = '  -3.58%  '
$ws.Range("E27").Value = '  -9.93%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.170'
$ws.Range("E28").Value = '  -5.77%  '
$ws.Range("E29").Value = '  -0.14%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '7.11'
$ws.Range("E30").Value = '  -7.48%  '
$ws.Range("E31").Value = '  +0.00%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '2.02'
$ws.Range("E32").Value = '  -4.71%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '7.46'
$ws.Range("E33").Value = '  -5.52%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '22.44'
$ws.Range("E34").Value = '  -3.34%  '
$ws.Range("E35").Value = '  -8.39%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '5.05'
$ws.Range("E36").Value = '  -7.48%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '166.04'
$ws.Range("E37").Value = '  -1.50%  '
$ws.Range("E38").Value = '  -6.49%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '6.63'
$ws.Range("E39").Value = '  -4.28%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '3.325.75'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '26.21'
$ws.Range("E41").Value = '  -15.73%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.0724'
$ws.Range("E42").Value = '  -7.66%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '41.73'
$ws.Range("E43").Value = '  -2.40%  '
$ws.Range("E44").Value = '  -3.86%  '
$ws.Range("E45").Value = '  -6.97%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.10'
$ws.Range("E46").Value = '  -5.90%  '
$ws.Range("E47").Value = '  -7.18%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.00'
$ws.Range("E48").Value = '  -0.07%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '2.337.01'
$ws.Range("E49").Value = '  -8.07%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '6.34'
$ws.Range("E50").Value = '  -7.68%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '21.21'
$ws.Range("E51").Value = '  -6.04%  '
